$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue $ws 'D2' '312.82'
Set-TextValue $ws 'E2' '1.01%'
Set-TextValue $ws 'G2' '3'
Set-TextValue $ws 'D3' '37.77'
Set-TextValue $ws 'E3' '-0.67%'
Set-TextValue $ws 'G3' '3'
Set-TextValue $ws 'D4' '5.137'
Set-TextValue $ws 'E4' '1.33%'
Set-TextValue $ws 'G4' '3'
Set-TextValue $ws 'D5' '0.07908'
Set-TextValue $ws 'E5' '1.71%'
Set-TextValue $ws 'G5' '3'
Set-TextValue $ws 'B6' 'GateToken'
Set-TextValue $ws 'C6' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D6' '4.416'
Set-TextValue $ws 'E6' '1.39%'
Set-TextValue $ws 'G6' '3'
Set-TextValue $ws 'D7' '1.909'
Set-TextValue $ws 'E7' '0.73%'
Set-TextValue $ws 'G7' '3'
Set-TextValue $ws 'B8' 'KuCoinToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws 'D8' '8.284'
Set-TextValue $ws 'E8' '1.13%'
Set-TextValue $ws 'G8' '3'
Set-TextValue $ws 'B9' 'BTSEToken'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws 'D9' '2.993'
Set-TextValue $ws 'E9' '-3.25%'
Set-TextValue $ws 'G9' '3'
Set-TextValue $ws 'B10' 'MXToken'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D10' '0.9235'
Set-TextValue $ws 'E10' '0.39%'
Set-TextValue $ws 'G10' '3'
Set-TextValue $ws 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D11' '0.1138'
Set-TextValue $ws 'E11' '-9.43%'
Set-TextValue $ws 'G11' '3'
Set-TextValue $ws 'B12' 'WazirX'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D12' '0.1903'
Set-TextValue $ws 'E12' '1.06%'
Set-TextValue $ws 'G12' '3'
Set-TextValue $ws 'B13' 'MandalaExchangeToken'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D13' '0.09114'
Set-TextValue $ws 'E13' '3.83%'
Set-TextValue $ws 'G13' '3'
Set-TextValue $ws 'B14' 'BitrueCoin'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D14' '0.03326'
Set-TextValue $ws 'E14' '-2.37%'
Set-TextValue $ws 'G14' '3'
Set-TextValue $ws 'B15' 'BitMartToken'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D15' '0.09620'
Set-TextValue $ws 'E15' '-0.95%'
Set-TextValue $ws 'G15' '3'
Set-TextValue $ws 'B16' 'BitForexToken'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D16' '0.001384'
Set-TextValue $ws 'E16' '0.95%'
Set-TextValue $ws 'G16' '3'
Set-TextValue $ws 'B17' 'TigerCash'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D17' '0.006092'
Set-TextValue $ws 'E17' '0.86%'
Set-TextValue $ws 'G17' '3'
Set-TextValue $ws 'B18' 'LEO'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D18' '3.563'
Set-TextValue $ws 'E18' '-0.34%'
Set-TextValue $ws 'G18' '3'
Set-TextValue $ws 'E19' '1.07%'
Set-TextValue $ws 'G19' '3'
Set-TextValue $ws 'D20' '5.897'
Set-TextValue $ws 'E20' '17.28%'
Set-TextValue $ws 'G20' '3'
Set-TextValue $ws 'D21' '0.1288'
Set-TextValue $ws 'E21' '0.69%'
Set-TextValue $ws 'G21' '3'
Set-TextValue $ws 'D22' '0.2592'
Set-TextValue $ws 'E22' '0.07%'
Set-TextValue $ws 'G22' '3'
Set-TextValue $ws 'D23' '0.04359'
Set-TextValue $ws 'E23' '-0.91%'
Set-TextValue $ws 'G23' '3'
Set-TextValue $ws 'D24' '0.001235'
Set-TextValue $ws 'E24' '1.97%'
Set-TextValue $ws 'G24' '3'
Set-TextValue $ws 'D25' '0.004644'
Set-TextValue $ws 'E25' '9.03%'
Set-TextValue $ws 'G25' '3'
Set-TextValue $ws 'D26' '0.0001362'
Set-TextValue $ws 'E26' '0.80%'
Set-TextValue $ws 'G26' '3'
Set-TextValue $ws 'D27' '0.0003995'
Set-TextValue $ws 'G27' '3'
Set-TextValue $ws 'G28' '3'
Set-TextValue $ws 'G29' '3'
Set-TextValue $ws 'G30' '3'
Set-TextValue $ws 'G31' '3'
Set-TextValue $ws 'G32' '3'
Set-TextValue $ws 'G33' '3'
Set-TextValue $ws 'G34' '3'
Set-TextValue $ws 'G35' '3'
Set-TextValue $ws 'G36' '3'
Set-TextValue $ws 'G37' '3'
Set-TextValue $ws 'G38' '3'
Set-TextValue $ws 'D39' '0.02273'
Set-TextValue $ws 'E39' '6.50%'
Set-TextValue $ws 'G39' '3'
Set-TextValue $ws 'D40' '0.05104'
Set-TextValue $ws 'E40' '2.30%'
Set-TextValue $ws 'G40' '3'
Set-TextValue $ws 'D41' '0.007470'
Set-TextValue $ws 'E41' '-7.78%'
Set-TextValue $ws 'G41' '3'
Set-TextValue $ws 'D42' '0.009038'
Set-TextValue $ws 'E42' '-10.46%'
Set-TextValue $ws 'G42' '3'
Set-TextValue $ws 'D43' '0.1357'
Set-TextValue $ws 'E43' '0.96%'
Set-TextValue $ws 'G43' '3'
Set-TextValue $ws 'D44' '0.001953'
Set-TextValue $ws 'E44' '-5.29%'
Set-TextValue $ws 'G44' '3'
Set-TextValue $ws 'D45' '0.008636'
Set-TextValue $ws 'E45' '-0.80%'
Set-TextValue $ws 'G45' '3'
Set-TextValue $ws 'D46' '0.00006683'
Set-TextValue $ws 'E46' '3.89%'
Set-TextValue $ws 'G46' '3'
Set-TextValue $ws 'G47' '3'
Set-TextValue $ws 'B48' 'BOLO'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws 'D48' '0.003284'
Set-TextValue $ws 'E48' '-3.23%'
Set-TextValue $ws 'G48' '3'
Set-TextValue $ws 'B49' 'CoinbaseStockToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws 'D49' '0.001001'
Set-TextValue $ws 'E49' '-40.71%'
Set-TextValue $ws 'G49' '3'
Set-TextValue $ws 'D50' '0.00002103'
Set-TextValue $ws 'G50' '3'
Set-TextValue $ws 'D51' '0.0002002'
Set-TextValue $ws 'G51' '3'
